$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Execucio" (B) and "BDD" (C) counts for week rows 2 and 3 ---
$ws.Range("B2").Value = 504
$ws.Range("C2").Value = 650
$ws.Range("B3").Value = 375
$ws.Range("C3").Value = 450

# --- New column E: ratio of BDD / Setmana for each row ---
$ws.Range("E2").Formula = "=C2/D2"
$ws.Range("E3").Formula = "=C3/D3"

# --- New row 4: week-over-week growth ratios for B, C, D ---
$ws.Range("B4").Formula = "=B3/B2"
$ws.Range("C4").Formula = "=C3/C2"
$ws.Range("D4").Formula = "=D3/D2"

# --- Formatting for the new percentage cells: centered, 0.00% ---
$ws.Range("E2:E3").NumberFormat = "0.00%"
$ws.Range("E2:E3").HorizontalAlignment = -4108
$ws.Range("B4:D4").NumberFormat = "0.00%"
$ws.Range("B4:D4").HorizontalAlignment = -4108

# --- Column width tweaks ---
$ws.Columns.Item(1).ColumnWidth = 3.9999999999999996
$ws.Columns.Item(2).ColumnWidth = 8.166666666666666
$ws.Columns.Item(3).ColumnWidth = 6.666666666666667
$ws.Columns.Item(4).ColumnWidth = 7.833333333333333
$ws.Columns.Item(5).ColumnWidth = 6.666666666666667

# --- Header/footer font style: "Normal" -> "Regular" ---
$ws.PageSetup.CenterHeader = "&""Times New Roman,Regular""&12&A"
$ws.PageSetup.CenterFooter = "&""Times New Roman,Regular""&12P" + [char]0x00E0 + "gina &P"
